# Generate Report for handoff
# Adds a new e2e test file "fa22870f-5f8d-45cb-9326-b2a18e178bd2" to the
# localization status report, and updates "e6359ff2-096d-43db-9062-628c8a4e76c7"
# from "Ready for handoff" to "In Translation" (with target/handback columns
# populated) on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Insert a new row above the "6f4488dc" row (currently row 5) to hold the
# refreshed "e6359ff2" entry; everything below shifts down by one.
$ws1.Rows.Item(5).Insert()

# New row 5: e6359ff2, now "In Translation"
$ws1.Range("B5").Value = "In Translation"
$ws1.Range("C5").Value = "In Translation"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0a97a9dccc5509b5324808d0f54e95002a0cfc1d/e2e/e6359ff2-096d-43db-9062-628c8a4e76c7.md", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.md")

# Row 7 now holds the stale "e6359ff2" data (shifted from old row 6); turn it
# into the brand new "fa22870f" entry instead.
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1b1c2d3e4f5061728394a5b6c7d8e9f0a1b2c3d4/e2e/fa22870f-5f8d-45cb-9326-b2a18e178bd2.md", "", "", "fa22870f-5f8d-45cb-9326-b2a18e178bd2.md")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(5).Insert()

# New row 5: e6359ff2, now "In Translation" with target/handback columns.
$ws2.Range("B5").Value = "In Translation"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0a97a9dccc5509b5324808d0f54e95002a0cfc1d/e2e/e6359ff2-096d-43db-9062-628c8a4e76c7.md", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f579a1bce1e106ca03dcb0d7d3ebf8ed830095c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.zh-cn.xlf", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.zh-cn.xlf")
$ws2.Range("D5").Value = "2016-01-19 04:06:52"
$ws2.Hyperlinks.Add($ws2.Range("E5"), "https://github.com/OpenLocalizationTest/oltest/blob/0a97a9dccc5509b5324808d0f54e95002a0cfc1d/e2e/e6359ff2-096d-43db-9062-628c8a4e76c7.md", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.md")
$ws2.Hyperlinks.Add($ws2.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f579a1bce1e106ca03dcb0d7d3ebf8ed830095c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.zh-cn.xlf", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.zh-cn.xlf")
$ws2.Range("G5").Value = "2016-01-19 04:07:35"
$ws2.Range("H5").Value = "Include"

# Row 7 now holds the stale "e6359ff2" data (shifted from old row 6); turn it
# into the brand new "fa22870f" entry instead.
$ws2.Range("B7").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1b1c2d3e4f5061728394a5b6c7d8e9f0a1b2c3d4/e2e/fa22870f-5f8d-45cb-9326-b2a18e178bd2.md", "", "", "fa22870f-5f8d-45cb-9326-b2a18e178bd2.md")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1ba51936f33433239ab2f5cd50f6dab1e511777/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fa22870f-5f8d-45cb-9326-b2a18e178bd2.c1ba51936f33433239ab2f5cd50f6dab1e511777.zh-cn.xlf", "", "", "fa22870f-5f8d-45cb-9326-b2a18e178bd2.c1ba51936f33433239ab2f5cd50f6dab1e511777.zh-cn.xlf")
$ws2.Range("D7").Value = "2016-01-19 04:09:04"
$ws2.Range("G7").Value = "0001-01-01 00:00:00"
$ws2.Range("H7").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(5).Insert()

# New row 5: e6359ff2, now "In Translation" with target/handback columns.
$ws3.Range("B5").Value = "In Translation"
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0a97a9dccc5509b5324808d0f54e95002a0cfc1d/e2e/e6359ff2-096d-43db-9062-628c8a4e76c7.md", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8484116cc101a660f93b7acae94923c8db39a05a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.de-de.xlf", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.de-de.xlf")
$ws3.Range("D5").Value = "2016-01-19 04:07:01"
$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://github.com/OpenLocalizationTest/oltest/blob/0a97a9dccc5509b5324808d0f54e95002a0cfc1d/e2e/e6359ff2-096d-43db-9062-628c8a4e76c7.md", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.md")
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8484116cc101a660f93b7acae94923c8db39a05a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.de-de.xlf", "", "", "e6359ff2-096d-43db-9062-628c8a4e76c7.4ebf6dfea1744e14015297824db9271adbfed053.de-de.xlf")
$ws3.Range("G5").Value = "2016-01-19 04:07:51"
$ws3.Range("H5").Value = "Include"

# Row 7 now holds the stale "e6359ff2" data (shifted from old row 6); turn it
# into the brand new "fa22870f" entry instead.
$ws3.Range("B7").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1b1c2d3e4f5061728394a5b6c7d8e9f0a1b2c3d4/e2e/fa22870f-5f8d-45cb-9326-b2a18e178bd2.md", "", "", "fa22870f-5f8d-45cb-9326-b2a18e178bd2.md")
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1ba51936f33433239ab2f5cd50f6dab1e511777/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fa22870f-5f8d-45cb-9326-b2a18e178bd2.c1ba51936f33433239ab2f5cd50f6dab1e511777.de-de.xlf", "", "", "fa22870f-5f8d-45cb-9326-b2a18e178bd2.c1ba51936f33433239ab2f5cd50f6dab1e511777.de-de.xlf")
$ws3.Range("D7").Value = "2016-01-19 04:09:14"
$ws3.Range("G7").Value = "0001-01-01 00:00:00"
$ws3.Range("H7").Value = "Include"

Write-Host "done"
